# Apply the "LinuxForHealth" rebrand / 8.0.0 release update to the
# StructureDefinition-insight-type workbook.
#
# Sheet "Metadata" (key/value table in columns A/B):
#   B2 - URL:       http://ibm.com/...            -> http://linuxforhealth.org/...
#   B3 - Version:   7.0.0                          -> 8.0.0
#   B8 - Date:      2022-09-08T16:11:15+00:00      -> 2022-11-10T16:00:46+00:00
#   B9 - Publisher: Alvearie Team                  -> LinuxForHealth Team
#
# Sheet "Elements" (generated FHIR StructureDefinition element table):
#   Q5  (Extension.url / Fixed Value) mirrors the same URL as Metadata!B2
#       and must be updated to match.
#   AI2 (Extension / Constraint(s)) previously held the ele-1/ext-1
#       constraint text; in the regenerated export that text now lives on
#       the Extension.extension row (AI4) instead, so AI2 is cleared.

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsMetadata.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/insight-type"
$wsMetadata.Range("B3").Value = "8.0.0"
$wsMetadata.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$wsMetadata.Range("B9").Value = "LinuxForHealth Team"

$wsElements = $wb.Worksheets.Item("Elements")
$wsElements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/insight-type"
$wsElements.Range("AI2").Value = ""
